# Updates cryptos list (prices, 1h volume %, and re-ordered/renamed coin rows)
# per commit: "Updated cryptos list on Sat Nov  4 23:45:22 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.451.30"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.877.89"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'240.33"
$ws.Range("E5").Value = "  +3.65%  "
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'42.88"
$ws.Range("E8").Value = "  +7.71%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "2.150.00"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "'11.62"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.873.33"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.686"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "'4.75"
$ws.Range("D17").Value = "35.421.70"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "'71.16"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "0.0₃0804"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").Value = "'243.33"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'170.44"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").Value = "'8.25"
$ws.Range("D27").Value = "'1.89"
$ws.Range("E27").Value = "  +24.28%  "
$ws.Range("D28").Value = "'17.85"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").Value = "'0.0565"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  +2.84%  "
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.06"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  +23.47%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.830"
$ws.Range("E35").Value = "  +18.67%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.05"
$ws.Range("E36").Value = "  +6.85%  "
$ws.Range("E37").Value = "  +7.94%  "
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("D39").Value = "'0.0204"
$ws.Range("E39").Value = "  +5.09%  "
$ws.Range("D40").Value = "'91.11"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").Value = "1.353.69"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.0605"
$ws.Range("E42").Value = "  +15.56%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'15.26"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("D44").Value = "'2.36"
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").Value = "'12.96"
$ws.Range("E45").Value = "  +56.40%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "'6.67"
$ws.Range("E47").Value = "  +6.70%  "
$ws.Range("D48").Value = "'2.73"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("D49").Value = "2.059.94"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").Value = "'0.0691"
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("E51").Value = "  +0.34%  "
